# gob_template.xlsx update
#  - "Resumen" summary sheet becomes the active/selected sheet (was "year")
#  - the placeholder tokens on the "Resumen" sheet are renamed from the
#    `results.*` namespace to `summary.*`, and a third placeholder
#    ({summary.titles}) is added in column C (column B, previously blank,
#    now carries {summary.graduates})
#  - column widths on "Resumen" are widened to fit the new/renamed labels
#    and row 2 grows taller to fit the wrapped text

$wb = $excel.ActiveWorkbook

$resumen = $wb.Worksheets.Item("Resumen")
$year = $wb.Worksheets.Item("year")

# --- Resumen!A2:C2 — rename + add the summary placeholders ----------------
# before: A2={results.year}  B2=(blank)            C2={results.count}
# after : A2={summary.year}  B2={summary.graduates} C2={summary.titles}
$resumen.Range("A2").Value = "{summary.year}"
$resumen.Range("B2").Value = "{summary.graduates}"
$resumen.Range("C2").Value = "{summary.titles}"

# Row 2 now wraps a longer label, so it grows taller.
$resumen.Rows.Item(2).RowHeight = 29.85

# --- column widths ----------------------------------------------------------
$resumen.Columns.Item(1).ColumnWidth = 15.6
$resumen.Columns.Item(2).ColumnWidth = 22.1
$resumen.Columns.Item(3).ColumnWidth = 24.43

# --- make "Resumen" the active sheet / selection ---------------------------
$resumen.Activate() | Out-Null
$resumen.Range("C3").Select() | Out-Null

# "year" keeps its own prior selection (F2); it is simply no longer the
# active tab once "Resumen" is activated above.
